$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Cells.Item(12, 8).Value = 444.4   # H12: 549.75 -> 444.4
$ws.Cells.Item(12, 10).Value = 86.5   # J12: 150 -> 86.5
$ws.Cells.Item(12, 12).Value = 86.5   # L12: 150 -> 86.5
$ws.Cells.Item(12, 14).Value = -426.5   # N12: -490 -> -426.5
# Row 51
$ws.Cells.Item(51, 8).Value = 10645.357   # H51: 9913.333000000001 -> 10645.357
$ws.Cells.Item(51, 9).Value = 9699.333000000001   # I51: 9813.333000000001 -> 9699.333000000001
$ws.Cells.Item(51, 10).Value = 10758.88   # J51: 9933.333000000001 -> 10758.88
$ws.Cells.Item(51, 11).Value = 9699.333000000001   # K51: 9813.333000000001 -> 9699.333000000001
$ws.Cells.Item(51, 12).Value = 10758.88   # L51: 9933.333000000001 -> 10758.88
$ws.Cells.Item(51, 13).Value = -9215.333000000001   # M51: -9329.333000000001 -> -9215.333000000001
$ws.Cells.Item(51, 14).Value = -11726.88   # N51: -10901.333 -> -11726.88
# Row 88
$ws.Cells.Item(88, 8).Value = 11719.692   # H88: 9081.091 -> 11719.692
$ws.Cells.Item(88, 9).Value = 1681.5   # I88: 2133 -> 1681.5
$ws.Cells.Item(88, 10).Value = 13544.818   # J88: 11686.625 -> 13544.818
$ws.Cells.Item(88, 11).Value = 1681.5   # K88: 2133 -> 1681.5
$ws.Cells.Item(88, 12).Value = 13544.818   # L88: 11686.625 -> 13544.818
$ws.Cells.Item(88, 13).Value = -1275.5   # M88: -1727 -> -1275.5
$ws.Cells.Item(88, 14).Value = -14356.818   # N88: -12498.625 -> -14356.818
# Row 91
$ws.Cells.Item(91, 8).Value = 11719.692   # H91: 9081.091 -> 11719.692
$ws.Cells.Item(91, 9).Value = 1681.5   # I91: 2133 -> 1681.5
$ws.Cells.Item(91, 10).Value = 13544.818   # J91: 11686.625 -> 13544.818
$ws.Cells.Item(91, 11).Value = 1681.5   # K91: 2133 -> 1681.5
$ws.Cells.Item(91, 12).Value = 13544.818   # L91: 11686.625 -> 13544.818
$ws.Cells.Item(91, 13).Value = -277.5   # M91: -729 -> -277.5
$ws.Cells.Item(91, 14).Value = -16352.818   # N91: -14494.625 -> -16352.818
# Row 112
$ws.Cells.Item(112, 8).Value = 4065.75   # H112: 3809.8518 -> 4065.75
$ws.Cells.Item(112, 10).Value = 4065.75   # J112: 3809.8518 -> 4065.75
$ws.Cells.Item(112, 12).Value = 12197.25   # L112: 11429.5554 -> 12197.25
$ws.Cells.Item(112, 14).Value = -14413.25   # N112: -13645.5554 -> -14413.25
# Row 125
$ws.Cells.Item(125, 8).Value = 6642.091   # H125: 6351.1816 -> 6642.091
$ws.Cells.Item(125, 9).Value = 815   # I125: 812.5 -> 815
$ws.Cells.Item(125, 10).Value = 11498   # J125: 12997.6 -> 11498
$ws.Cells.Item(125, 11).Value = 7335   # K125: 7312.5 -> 7335
$ws.Cells.Item(125, 12).Value = 103482   # L125: 116978.4 -> 103482
$ws.Cells.Item(125, 13).Value = -4875   # M125: -4852.5 -> -4875
$ws.Cells.Item(125, 14).Value = -108402   # N125: -121898.4 -> -108402
# Row 132
$ws.Cells.Item(132, 8).Value = 4796.4   # H132: 5112.6787 -> 4796.4
$ws.Cells.Item(132, 9).Value = 1585.7916   # I132: 1638.7391 -> 1585.7916
$ws.Cells.Item(132, 10).Value = 17638.834   # J132: 21092.8 -> 17638.834
$ws.Cells.Item(132, 11).Value = 4757.3748   # K132: 4916.2173 -> 4757.3748
$ws.Cells.Item(132, 12).Value = 52916.50199999999   # L132: 63278.39999999999 -> 52916.50199999999
$ws.Cells.Item(132, 13).Value = -2227.3748   # M132: -2386.2173 -> -2227.3748
$ws.Cells.Item(132, 14).Value = -57976.50199999999   # N132: -68338.39999999999 -> -57976.50199999999
# Row 138
$ws.Cells.Item(138, 8).Value = 5315.4336   # H138: 5547.5923 -> 5315.4336
$ws.Cells.Item(138, 9).Value = 946.0909   # I138: 1020.4545 -> 946.0909
$ws.Cells.Item(138, 10).Value = 6891.262   # J138: 7391.9814 -> 6891.262
$ws.Cells.Item(138, 11).Value = 2838.2727   # K138: 3061.3635 -> 2838.2727
$ws.Cells.Item(138, 12).Value = 20673.786   # L138: 22175.9442 -> 20673.786
$ws.Cells.Item(138, 13).Value = 2301.7273   # M138: 2078.6365 -> 2301.7273
$ws.Cells.Item(138, 14).Value = -30953.786   # N138: -32455.9442 -> -30953.786
# Row 140
$ws.Cells.Item(140, 8).Value = 71507.17999999999   # H140: 71319 -> 71507.17999999999
$ws.Cells.Item(140, 10).Value = 71587   # J140: 71374.45 -> 71587
$ws.Cells.Item(140, 12).Value = 71587   # L140: 71374.45 -> 71587
$ws.Cells.Item(140, 14).Value = -81947   # N140: -81734.45 -> -81947

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 18224.229   # H32: 18079.404 -> 18224.229
$ws.Cells.Item(32, 9).Value = 17494.682   # I32: 17493.977 -> 17494.682
$ws.Cells.Item(32, 10).Value = 26249.25   # J32: 26665.666 -> 26249.25
$ws.Cells.Item(32, 11).Value = 17494.682   # K32: 17493.977 -> 17494.682
$ws.Cells.Item(32, 12).Value = 26249.25   # L32: 26665.666 -> 26249.25
$ws.Cells.Item(32, 13).Value = -17207.682   # M32: -17206.977 -> -17207.682
$ws.Cells.Item(32, 14).Value = -26823.25   # N32: -27239.666 -> -26823.25
# Row 61
$ws.Cells.Item(61, 8).Value = 2331.375   # H61: 2407.4285 -> 2331.375
$ws.Cells.Item(61, 10).Value = 1799   # J61: 0 -> 1799
$ws.Cells.Item(61, 12).Value = 1799   # L61: 0 -> 1799
$ws.Cells.Item(61, 14).Value = -2223   # N61: None -> -2223
# Row 74
$ws.Cells.Item(74, 8).Value = 1181.5264   # H74: 1181.7894 -> 1181.5264
$ws.Cells.Item(74, 9).Value = 1159.1177   # I74: 1136.3334 -> 1159.1177
$ws.Cells.Item(74, 10).Value = 1372   # J74: 2000 -> 1372
$ws.Cells.Item(74, 11).Value = 1159.1177   # K74: 1136.3334 -> 1159.1177
$ws.Cells.Item(74, 12).Value = 1372   # L74: 2000 -> 1372
$ws.Cells.Item(74, 13).Value = -285.1177   # M74: -262.3334 -> -285.1177
$ws.Cells.Item(74, 14).Value = -3120   # N74: -3748 -> -3120
# Row 77
$ws.Cells.Item(77, 8).Value = 1181.5264   # H77: 1181.7894 -> 1181.5264
$ws.Cells.Item(77, 9).Value = 1159.1177   # I77: 1136.3334 -> 1159.1177
$ws.Cells.Item(77, 10).Value = 1372   # J77: 2000 -> 1372
$ws.Cells.Item(77, 11).Value = 5795.5885   # K77: 5681.666999999999 -> 5795.5885
$ws.Cells.Item(77, 12).Value = 6860   # L77: 10000 -> 6860
$ws.Cells.Item(77, 13).Value = -1427.5885   # M77: -1313.666999999999 -> -1427.5885
$ws.Cells.Item(77, 14).Value = -15596   # N77: -18736 -> -15596
# Row 132
$ws.Cells.Item(132, 8).Value = 1825   # H132: 1953.7333 -> 1825
$ws.Cells.Item(132, 9).Value = 1825   # I132: 1953.7333 -> 1825
$ws.Cells.Item(132, 11).Value = 5475   # K132: 5861.199900000001 -> 5475
$ws.Cells.Item(132, 13).Value = -2945   # M132: -3331.199900000001 -> -2945
# Row 136
$ws.Cells.Item(136, 8).Value = 2331.375   # H136: 2407.4285 -> 2331.375
$ws.Cells.Item(136, 10).Value = 1799   # J136: 0 -> 1799
$ws.Cells.Item(136, 12).Value = 5397   # L136: 0 -> 5397
$ws.Cells.Item(136, 14).Value = -10497   # N136: None -> -10497

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Cells.Item(134, 8).Value = 2863.697   # H134: 2909.5625 -> 2863.697
$ws.Cells.Item(134, 9).Value = 2184.6316   # I134: 2228.4443 -> 2184.6316
$ws.Cells.Item(134, 11).Value = 6553.8948   # K134: 6685.3329 -> 6553.8948
$ws.Cells.Item(134, 13).Value = -4018.8948   # M134: -4150.3329 -> -4018.8948

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 115
$ws.Cells.Item(115, 8).Value = 49665.668   # H115: 49247.332 -> 49665.668
$ws.Cells.Item(115, 10).Value = 49665.668   # J115: 49247.332 -> 49665.668
$ws.Cells.Item(115, 12).Value = 49665.668   # L115: 49247.332 -> 49665.668
$ws.Cells.Item(115, 14).Value = -52015.668   # N115: -51597.332 -> -52015.668
# Row 132
$ws.Cells.Item(132, 8).Value = 102579050   # H132: 121229140 -> 102579050
$ws.Cells.Item(132, 9).Value = 190485380   # I132: 222232590 -> 190485380
$ws.Cells.Item(132, 10).Value = 21666.334   # J132: 24999.8 -> 21666.334
$ws.Cells.Item(132, 11).Value = 571456140   # K132: 666697770 -> 571456140
$ws.Cells.Item(132, 12).Value = 64999.00199999999   # L132: 74999.39999999999 -> 64999.00199999999
$ws.Cells.Item(132, 13).Value = -571453610   # M132: -666695240 -> -571453610
$ws.Cells.Item(132, 14).Value = -70059.00199999999   # N132: -80059.39999999999 -> -70059.00199999999
# Row 133
$ws.Cells.Item(133, 8).Value = 59422.855   # H133: 54550 -> 59422.855
$ws.Cells.Item(133, 9).Value = 0   # I133: 25000 -> 0
$ws.Cells.Item(133, 10).Value = 59422.855   # J133: 58243.75 -> 59422.855
$ws.Cells.Item(133, 11).Value = 0   # K133: 25000 -> 0
$ws.Cells.Item(133, 12).Value = 59422.855   # L133: 58243.75 -> 59422.855
$ws.Cells.Item(133, 13).ClearContents()   # M133: delete (was -22470)
$ws.Cells.Item(133, 14).Value = -64482.855   # N133: -63303.75 -> -64482.855
# Row 134
$ws.Cells.Item(134, 8).Value = 2845.4443   # H134: 3062.7693 -> 2845.4443
$ws.Cells.Item(134, 9).Value = 2190.3333   # I134: 2373.182 -> 2190.3333
$ws.Cells.Item(134, 10).Value = 3369.5334   # J134: 3568.4666 -> 3369.5334
$ws.Cells.Item(134, 11).Value = 6570.999899999999   # K134: 7119.545999999999 -> 6570.999899999999
$ws.Cells.Item(134, 12).Value = 10108.6002   # L134: 10705.3998 -> 10108.6002
$ws.Cells.Item(134, 13).Value = -4035.999899999999   # M134: -4584.545999999999 -> -4035.999899999999
$ws.Cells.Item(134, 14).Value = -15178.6002   # N134: -15775.3998 -> -15178.6002

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 644.7059   # H5: 616.5 -> 644.7059
$ws.Cells.Item(5, 10).Value = 1258   # J5: 1298.5 -> 1258
$ws.Cells.Item(5, 12).Value = 3774   # L5: 3895.5 -> 3774
$ws.Cells.Item(5, 14).Value = -3998   # N5: -4119.5 -> -3998
# Row 14
$ws.Cells.Item(14, 8).Value = 213.44444   # H14: 237.94444 -> 213.44444
$ws.Cells.Item(14, 9).Value = 213.44444   # I14: 237.94444 -> 213.44444
$ws.Cells.Item(14, 11).Value = 640.33332   # K14: 713.83332 -> 640.33332
$ws.Cells.Item(14, 13).Value = -467.33332   # M14: -540.83332 -> -467.33332
# Row 37
$ws.Cells.Item(37, 8).Value = 111210890   # H37: 100106980 -> 111210890
$ws.Cells.Item(37, 10).Value = 111210890   # J37: 100106980 -> 111210890
$ws.Cells.Item(37, 12).Value = 333632670   # L37: 300320940 -> 333632670
$ws.Cells.Item(37, 14).Value = -333632894   # N37: -300321164 -> -333632894
# Row 107
$ws.Cells.Item(107, 8).Value = 873.1724   # H107: 920.0741 -> 873.1724
$ws.Cells.Item(107, 9).Value = 720.9   # I107: 841.125 -> 720.9
$ws.Cells.Item(107, 11).Value = 2162.7   # K107: 2523.375 -> 2162.7
$ws.Cells.Item(107, 13).Value = -242.6999999999998   # M107: -603.375 -> -242.6999999999998
# Row 114
$ws.Cells.Item(114, 8).Value = 4065.5   # H114: 5598.5 -> 4065.5
$ws.Cells.Item(114, 10).Value = 4065.5   # J114: 5598.5 -> 4065.5
$ws.Cells.Item(114, 12).Value = 12196.5   # L114: 16795.5 -> 12196.5
$ws.Cells.Item(114, 14).Value = -18704.5   # N114: -23303.5 -> -18704.5
# Row 132
$ws.Cells.Item(132, 8).Value = 1233.5758   # H132: 1187.6875 -> 1233.5758
$ws.Cells.Item(132, 9).Value = 1024.12   # I132: 981.5185 -> 1024.12
$ws.Cells.Item(132, 10).Value = 1888.125   # J132: 2301 -> 1888.125
$ws.Cells.Item(132, 11).Value = 9217.079999999998   # K132: 8833.666499999999 -> 9217.079999999998
$ws.Cells.Item(132, 12).Value = 16993.125   # L132: 20709 -> 16993.125
$ws.Cells.Item(132, 13).Value = -6687.079999999998   # M132: -6303.666499999999 -> -6687.079999999998
$ws.Cells.Item(132, 14).Value = -22053.125   # N132: -25769 -> -22053.125
# Row 135
$ws.Cells.Item(135, 8).Value = 644.7059   # H135: 616.5 -> 644.7059
$ws.Cells.Item(135, 10).Value = 1258   # J135: 1298.5 -> 1258
$ws.Cells.Item(135, 12).Value = 11322   # L135: 11686.5 -> 11322
$ws.Cells.Item(135, 14).Value = -16392   # N135: -16756.5 -> -16392

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70, 8).Value = 8603.200000000001   # H70: 8631.666999999999 -> 8603.200000000001
$ws.Cells.Item(70, 9).Value = 8428.23   # I70: 8447.666999999999 -> 8428.23
$ws.Cells.Item(70, 10).Value = 8928.143   # J70: 8999.666999999999 -> 8928.143
$ws.Cells.Item(70, 11).Value = 8428.23   # K70: 8447.666999999999 -> 8428.23
$ws.Cells.Item(70, 12).Value = 8928.143   # L70: 8999.666999999999 -> 8928.143
$ws.Cells.Item(70, 13).Value = -8158.23   # M70: -8177.666999999999 -> -8158.23
$ws.Cells.Item(70, 14).Value = -9468.143   # N70: -9539.666999999999 -> -9468.143
# Row 73
$ws.Cells.Item(73, 8).Value = 8603.200000000001   # H73: 8631.666999999999 -> 8603.200000000001
$ws.Cells.Item(73, 9).Value = 8428.23   # I73: 8447.666999999999 -> 8428.23
$ws.Cells.Item(73, 10).Value = 8928.143   # J73: 8999.666999999999 -> 8928.143
$ws.Cells.Item(73, 11).Value = 8428.23   # K73: 8447.666999999999 -> 8428.23
$ws.Cells.Item(73, 12).Value = 8928.143   # L73: 8999.666999999999 -> 8928.143
$ws.Cells.Item(73, 13).Value = -7492.23   # M73: -7511.666999999999 -> -7492.23
$ws.Cells.Item(73, 14).Value = -10800.143   # N73: -10871.667 -> -10800.143
# Row 80
$ws.Cells.Item(80, 8).Value = 47491.19   # H80: 55735.047 -> 47491.19
$ws.Cells.Item(80, 9).Value = 55919.58   # I80: 70258.13 -> 55919.58
$ws.Cells.Item(80, 11).Value = 55919.58   # K80: 70258.13 -> 55919.58
$ws.Cells.Item(80, 13).Value = -54921.58   # M80: -69260.13 -> -54921.58
# Row 83
$ws.Cells.Item(83, 8).Value = 47491.19   # H83: 55735.047 -> 47491.19
$ws.Cells.Item(83, 9).Value = 55919.58   # I83: 70258.13 -> 55919.58
$ws.Cells.Item(83, 11).Value = 279597.9   # K83: 351290.65 -> 279597.9
$ws.Cells.Item(83, 13).Value = -274605.9   # M83: -346298.65 -> -274605.9

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Cells.Item(16, 8).Value = 1977.8889   # H16: 899.5 -> 1977.8889
$ws.Cells.Item(16, 9).Value = 2162.625   # I16: 899 -> 2162.625
$ws.Cells.Item(16, 10).Value = 500   # J16: 900 -> 500
$ws.Cells.Item(16, 11).Value = 2162.625   # K16: 899 -> 2162.625
$ws.Cells.Item(16, 12).Value = 500   # L16: 900 -> 500
$ws.Cells.Item(16, 13).Value = -1992.625   # M16: -729 -> -1992.625
$ws.Cells.Item(16, 14).Value = -840   # N16: -1240 -> -840
# Row 22
$ws.Cells.Item(22, 8).Value = 706.1667   # H22: 720.7143 -> 706.1667
$ws.Cells.Item(22, 9).Value = 766.4286   # I22: 794.7692 -> 766.4286
$ws.Cells.Item(22, 10).Value = 667.8182   # J22: 676.9545000000001 -> 667.8182
$ws.Cells.Item(22, 11).Value = 766.4286   # K22: 794.7692 -> 766.4286
$ws.Cells.Item(22, 12).Value = 667.8182   # L22: 676.9545000000001 -> 667.8182
$ws.Cells.Item(22, 13).Value = -471.4286   # M22: -499.7692 -> -471.4286
$ws.Cells.Item(22, 14).Value = -1257.8182   # N22: -1266.9545 -> -1257.8182
# Row 27
$ws.Cells.Item(27, 8).Value = 706.1667   # H27: 720.7143 -> 706.1667
$ws.Cells.Item(27, 9).Value = 766.4286   # I27: 794.7692 -> 766.4286
$ws.Cells.Item(27, 10).Value = 667.8182   # J27: 676.9545000000001 -> 667.8182
$ws.Cells.Item(27, 11).Value = 766.4286   # K27: 794.7692 -> 766.4286
$ws.Cells.Item(27, 12).Value = 667.8182   # L27: 676.9545000000001 -> 667.8182
$ws.Cells.Item(27, 13).Value = -659.4286   # M27: -687.7692 -> -659.4286
$ws.Cells.Item(27, 14).Value = -881.8182   # N27: -890.9545000000001 -> -881.8182
# Row 129
$ws.Cells.Item(129, 8).Value = 149999   # H129: 149809 -> 149999
$ws.Cells.Item(129, 10).Value = 149999   # J129: 149809 -> 149999
$ws.Cells.Item(129, 12).Value = 149999   # L129: 149809 -> 149999
$ws.Cells.Item(129, 14).Value = -159999   # N129: -159809 -> -159999
# Row 132
$ws.Cells.Item(132, 8).Value = 4045.1765   # H132: 4063.602 -> 4045.1765
$ws.Cells.Item(132, 9).Value = 3247.0205   # I132: 3347.8215 -> 3247.0205
$ws.Cells.Item(132, 10).Value = 5131.5557   # J132: 5146.946 -> 5131.5557
$ws.Cells.Item(132, 11).Value = 9741.0615   # K132: 10043.4645 -> 9741.0615
$ws.Cells.Item(132, 12).Value = 15394.6671   # L132: 15440.838 -> 15394.6671
$ws.Cells.Item(132, 13).Value = -7211.0615   # M132: -7513.4645 -> -7211.0615
$ws.Cells.Item(132, 14).Value = -20454.6671   # N132: -20500.838 -> -20454.6671
# Row 136
$ws.Cells.Item(136, 8).Value = 4305.647   # H136: 3902 -> 4305.647
$ws.Cells.Item(136, 9).Value = 2834.5518   # I136: 2690.825 -> 2834.5518
$ws.Cells.Item(136, 10).Value = 6244.8184   # J136: 6209 -> 6244.8184
$ws.Cells.Item(136, 11).Value = 8503.6554   # K136: 8072.474999999999 -> 8503.6554
$ws.Cells.Item(136, 12).Value = 18734.4552   # L136: 18627 -> 18734.4552
$ws.Cells.Item(136, 13).Value = -5953.6554   # M136: -5522.474999999999 -> -5953.6554
$ws.Cells.Item(136, 14).Value = -23834.4552   # N136: -23727 -> -23834.4552

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Cells.Item(122, 8).Value = 5500.5   # H122: 5778.3335 -> 5500.5
$ws.Cells.Item(122, 9).Value = 3001.5   # I122: 3002 -> 3001.5
$ws.Cells.Item(122, 11).Value = 9004.5   # K122: 9006 -> 9004.5
$ws.Cells.Item(122, 13).Value = -6554.5   # M122: -6556 -> -6554.5
# Row 126
$ws.Cells.Item(126, 8).Value = 3005   # H126: 3004 -> 3005
$ws.Cells.Item(126, 9).Value = 0   # I126: 3003 -> 0
$ws.Cells.Item(126, 11).Value = 0   # K126: 9009 -> 0
$ws.Cells.Item(126, 13).ClearContents()   # M126: delete (was -6539)
# Row 136
$ws.Cells.Item(136, 8).Value = 3353.3416   # H136: 3181.6365 -> 3353.3416
$ws.Cells.Item(136, 9).Value = 2333.182   # I136: 2208.3333 -> 2333.182
$ws.Cells.Item(136, 11).Value = 6999.545999999999   # K136: 6624.999899999999 -> 6999.545999999999
$ws.Cells.Item(136, 13).Value = -4449.545999999999   # M136: -4074.999899999999 -> -4449.545999999999
